$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit reorders (permutes) the data rows of the trait/cell_type Fisher-test
# table (rows 3-53) without changing which row-numbers exist; the full 18-column
# (A:R) content of certain rows is relocated to other row numbers. We snapshot the
# full row ranges first, then write them back in their new positions so that no
# data is lost while rows are rotated between their old and new locations.

# Cycle covering rows: 3, 4
$row3 = $ws.Range("A3:R3").Value2
$row4 = $ws.Range("A4:R4").Value2
$ws.Range("A3:R3").Value2 = $row4
$ws.Range("A4:R4").Value2 = $row3

# Cycle covering rows: 10, 12
$row10 = $ws.Range("A10:R10").Value2
$row12 = $ws.Range("A12:R12").Value2
$ws.Range("A10:R10").Value2 = $row12
$ws.Range("A12:R12").Value2 = $row10

# Cycle covering rows: 17, 18
$row17 = $ws.Range("A17:R17").Value2
$row18 = $ws.Range("A18:R18").Value2
$ws.Range("A17:R17").Value2 = $row18
$ws.Range("A18:R18").Value2 = $row17

# Cycle covering rows: 19, 22, 20, 21
$row19 = $ws.Range("A19:R19").Value2
$row22 = $ws.Range("A22:R22").Value2
$row20 = $ws.Range("A20:R20").Value2
$row21 = $ws.Range("A21:R21").Value2
$ws.Range("A19:R19").Value2 = $row22
$ws.Range("A22:R22").Value2 = $row20
$ws.Range("A20:R20").Value2 = $row21
$ws.Range("A21:R21").Value2 = $row19

# Cycle covering rows: 36, 37
$row36 = $ws.Range("A36:R36").Value2
$row37 = $ws.Range("A37:R37").Value2
$ws.Range("A36:R36").Value2 = $row37
$ws.Range("A37:R37").Value2 = $row36

# Cycle covering rows: 39, 53, 40, 42, 51
$row39 = $ws.Range("A39:R39").Value2
$row53 = $ws.Range("A53:R53").Value2
$row40 = $ws.Range("A40:R40").Value2
$row42 = $ws.Range("A42:R42").Value2
$row51 = $ws.Range("A51:R51").Value2
$ws.Range("A39:R39").Value2 = $row53
$ws.Range("A53:R53").Value2 = $row40
$ws.Range("A40:R40").Value2 = $row42
$ws.Range("A42:R42").Value2 = $row51
$ws.Range("A51:R51").Value2 = $row39

# Cycle covering rows: 41, 47, 44, 48, 50, 45, 49
$row41 = $ws.Range("A41:R41").Value2
$row47 = $ws.Range("A47:R47").Value2
$row44 = $ws.Range("A44:R44").Value2
$row48 = $ws.Range("A48:R48").Value2
$row50 = $ws.Range("A50:R50").Value2
$row45 = $ws.Range("A45:R45").Value2
$row49 = $ws.Range("A49:R49").Value2
$ws.Range("A41:R41").Value2 = $row47
$ws.Range("A47:R47").Value2 = $row44
$ws.Range("A44:R44").Value2 = $row48
$ws.Range("A48:R48").Value2 = $row50
$ws.Range("A50:R50").Value2 = $row45
$ws.Range("A45:R45").Value2 = $row49
$ws.Range("A49:R49").Value2 = $row41

# Cycle covering rows: 43, 52, 46
$row43 = $ws.Range("A43:R43").Value2
$row52 = $ws.Range("A52:R52").Value2
$row46 = $ws.Range("A46:R46").Value2
$ws.Range("A43:R43").Value2 = $row52
$ws.Range("A52:R52").Value2 = $row46
$ws.Range("A46:R46").Value2 = $row43
